$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected cells are formatted as text so that values like "1.003"
# or "30.845.65" are preserved exactly as strings rather than being reinterpreted
# as numbers or dates by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.845.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.920.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4911"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2967"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06764"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.910.83"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.07"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07305"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6716"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.830.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007993"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.174.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.215"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "205.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.278"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.63%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.977"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.432"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.364"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09186"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05189"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7515"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.124"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.727"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01858"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.731"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9261"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.082"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4499"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.91"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +26.07%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.936"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.88%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1394"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.695"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.04"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.078"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05948"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4079"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.43%  "
